# Generate Report for Handoff
# - Update the "Latest HO Xliff Generate Date" (Overview!G) for the rows whose
#   handoff just completed.
# - Update the "Latest Handoff Datetime" (zh-cn!H / de-de!H) to the new handoff
#   time (de-de's handoff timestamp matches the Overview generate date; zh-cn's
#   is its own distinct timestamp).
# - Mark those rows' Priority (column E) as "ht" (handoff type) on both the
#   zh-cn and de-de localization sheets.
# Row 11 (894ba8db-...) is untouched, matching the source diff.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 10, 12, 13, 14)

# --- Overview sheet: bump the "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-17 16:21:37"
}

# --- zh-cn sheet: mark Priority as "ht" and refresh its own handoff datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-17 16:21:30"
}

# --- de-de sheet: mark Priority as "ht" and refresh its handoff datetime
#     (shares the same timestamp string as the Overview generate date) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-17 16:21:37"
}
